$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the content of row 3 and row 4 for columns
# A, B, E, F, G, H, P, Q, R, S, AC (other columns are identical
# between the two rows already, so no change is needed there).

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S", "AC")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"
    $val3 = $ws.Range($addr3).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value2 = $val4
    $ws.Range($addr4).Value2 = $val3
}
